# Update "New Hashed Password" column (C) for rows 2-26 on the
# "Password Audit" worksheet with newly generated bcrypt hashes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Password Audit")

$newHashes = @{
    2  = '$2b$10$g3pe90i7yVrJneH90WOn.O1FZCPS./1aXwd.2d1RsB1eIHGHIlOLG'
    3  = '$2b$10$mT/5RgxVZypYOBaRKina6OTK9NsyjZLJK8mHjMo26yhKiDayAwAqi'
    4  = '$2b$10$pMZqe8RF7Y4H4tG7rY6Tf.5VQGIPyAEM.ooqE6CecRknahO.Za3SC'
    5  = '$2b$10$XrDka72pDjPMXX6NEXx6a.uhr/5SQ1r/pyN5omWvUTM7GJENM7VAG'
    6  = '$2b$10$eUwCZTqheO0Vqp6C6P2LJuV7n4tc0UbkRvD04DQZjYJRZESbGcAQ6'
    7  = '$2b$10$vL8q2ZHPIGcqDczHXaz4YOrqF7ZLxKCQ6JreHjd7zB9yuUcqKVQSq'
    8  = '$2b$10$5z34zPmPRdZPGUkRxSW/KOybRVrmp6Zj1LsSUaPiMFyaHCIXTvUom'
    9  = '$2b$10$NrTJqztu/UtsqPQsKdBfWuGMNlPuCZssYMYHLlsjbWlBdlzIlVIRm'
    10 = '$2b$10$PKZGj4UrNjYU9hwj7oGOF.QpayoBANS3kxp9YSdLkMnBoK8b.s2O2'
    11 = '$2b$10$HoEU.ZazQNoZOf0n98m2A.7nDIO9bAFu9cU3kxwhJy5CNw87ml1rG'
    12 = '$2b$10$9ALv.H2NmoWaXzWT6F5vJuieLbnuSyKESdH8RQegNDv.v/qn1xpWu'
    13 = '$2b$10$sBPE0PFJPs4eTDa4ckta0OaOc62wrcfKmQAxyc30WTyVllDCNrhcm'
    14 = '$2b$10$JRkzFvSLfKI5hJ3EknyV7.o./Qa8DwkeuOqCTTtGo4UBPgTdK0X06'
    15 = '$2b$10$OAGgl7pREzHxi33EVzMKe.tatRR/eMKEHFEDFH.FHdx3HZZrjZ6fO'
    16 = '$2b$10$Vk658euGqnyNZj/80SXOleuICWZZo1JX1QOVwICr6HI46VCydtLmG'
    17 = '$2b$10$okyuJCgeUS/zUqL1lDGJvOjeGMQ6Q1D.zak4rzO7Yrdjddqi87Nru'
    18 = '$2b$10$mGs65vl5pjU0gjzTnbVC9u21.yOagrsB0p8BCYkg83ACNK397tpga'
    19 = '$2b$10$TwzC955E72vyJx0VCyYgG.WcIfEgqq22xcAkAQOsrF4IwVuR5SAvy'
    20 = '$2b$10$ViBzSWs9pVaKGPSbYBkOZOQRjd5m758q7YJ0U9fgbE/W3i59k1SyW'
    21 = '$2b$10$zbQos.imaLJhVepjErMaAerC99GUxFxOV2r34wjhpuOEEIGVahO0i'
    22 = '$2b$10$Fd1XSabXEwjiJNLocx6/J.gwwUC.I3RFgmMSnm8QZVB2kCvvX9DMS'
    23 = '$2b$10$nT5N8YyDvyfF6jBqdMrASeOppYMXTiDavONk6koDR.0prmDow3Zx2'
    24 = '$2b$10$HwahCWjiPns4h6MkeohBe.YQmvy7o66K/IeLlvIeiSPM5hr1zrmWq'
    25 = '$2b$10$UoiN1Gw07cn.16oD7YIww.qXM7HUrvegCndaNkNk5SuRGc0nQOHAS'
    26 = '$2b$10$JC9UU0nbTprkUJzlsgezf.V6knLGVWkSxevXt8CK3IwyLVbnTFmme'
}

foreach ($row in $newHashes.Keys) {
    $ws.Cells.Item($row, 3).Value = $newHashes[$row]
}
